$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 25 (CKD Stage): 1.0 -> 0.0, plus new C/D values (previously row 28's values)
# Force column B to stay text (it holds numeric-looking labels like "0.0","1.0",...)
$ws.Range("B25:B28").NumberFormat = "@"

$ws.Range("B25").Value = "0.0"
$ws.Range("C25").Value = "583 (91.5)"
$ws.Range("D25").Value = "2400 (93.1)"

$ws.Range("B26").Value = "1.0"
$ws.Range("C26").Value = "2 (0.3)"
$ws.Range("D26").Value = "1 (0.0)"

$ws.Range("B27").Value = "2.0"
$ws.Range("C27").Value = "8 (1.3)"
$ws.Range("D27").Value = "25 (1.0)"

$ws.Range("B28").Value = "3.0"
$ws.Range("C28").Value = "44 (6.9)"
$ws.Range("D28").Value = "153 (5.9)"
